$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56 (pushes old rows 56..194 down to 57..195)
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with a copy of the (now shifted) row 57's
# original data, but with an updated Fecha (date) and Volumen value.
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("D56").Value = 44526
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 100112037
$ws.Range("G56").Value = "Cebollín"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 180
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 6000
$ws.Range("M56").Value = 6000
$ws.Range("N56").Value = "`$/paquete 36 unidades"
$ws.Range("O56").Value = "Región Metropolitana"
$ws.Range("P56").Value = 167
$ws.Range("Q56").Value = 36
$ws.Range("R56").Value = "Hortaliza"
